$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "68D23006A"
$ws.Range("B2").Value = "[]"
